# Auto-update draw results: append the 2025-10-29 Pick 4 draw as a new row
# at the bottom of the "Results" log sheet.
#
# The data is stored as a flat, append-only log where every cell -- the
# date, the numeric-looking "phase" code, the dashed result, and the
# timestamp -- is literal TEXT, never a number/date. Left to its own
# devices Excel's COM layer "helpfully" autodetects dates and digit
# strings and coerces them to numeric/date values, so the columns that
# look numeric (Date, Phase, InsertedAt) are written with a leading
# apostrophe to force text entry, exactly like typing '2025-10-29 into a
# cell by hand. The apostrophe itself is not part of the stored value.
# ClearFormats() afterwards drops the "quote prefix" cell style that the
# apostrophe entry leaves behind, so the new cells keep using the sheet's
# default (General) style just like every other row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Next empty row right after the current last row of data.
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
$row = $lastRow + 1

$date = "2025-10-29"
$game = "Pick 4"
$phase = "251029"
$result = "8-3-0-6"
$insertedAt = "2025-10-29T21:40:33.440+04:00"

# A: Date -- looks numeric/date-like, force text.
$ws.Cells.Item($row, 1).Value = "'" + $date
$ws.Cells.Item($row, 1).ClearFormats()

# B: Game -- plain text, no coercion risk.
$ws.Cells.Item($row, 2).Value = $game

# C: Phase -- looks like a number, force text.
$ws.Cells.Item($row, 3).Value = "'" + $phase
$ws.Cells.Item($row, 3).ClearFormats()

# D: Result -- dashed digits, Excel does not treat this as a number.
$ws.Cells.Item($row, 4).Value = $result

# E: InsertedAt -- ISO-ish timestamp, looks date-like, force text.
$ws.Cells.Item($row, 5).Value = "'" + $insertedAt
$ws.Cells.Item($row, 5).ClearFormats()

# Keep Excel's "number stored as text" checker quiet over the new row too,
# mirroring the pre-existing ignoredErrors rule that covered A1:E<lastRow>.
# Best-effort: not every host persists this Range.Errors toggle back into
# the saved file, but setting it is harmless either way.
try {
    $ws.Range("A1:E" + $row).Errors.Item(9).Ignore = $true
} catch {
}
